$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight the "green" feature columns (Volatile acidity, Citric acid,
#     Total sulfur dioxide, alcohol) across the header + 2 model rows ---
$greenFill = 9818287      # RGB(175,208,149) in BGR-packed long used by Excel.Color
$ws.Range("C1:D3").Interior.Color = $greenFill
$ws.Range("G1:G3").Interior.Color = $greenFill
$ws.Range("L1:L3").Interior.Color = $greenFill

# --- New section: Linear regression model ---
$ws.Range("A7").Value = "Linear regression"

$ws.Range("A8").Value = "P-values"
$ws.Range("B8").Value = 0.470242504
$ws.Range("C8").Value = 0.0558321578
$ws.Range("D8").Value = 0.988234525
$ws.Range("E8").Value = [double]"1.24599779E-84"
$ws.Range("F8").Value = 0.908072679
$ws.Range("G8").Value = [double]"2.96703179E-181"
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0.999999999
$ws.Range("J8").Value = 0.999119348
$ws.Range("K8").Value = 0.997358586
$ws.Range("L8").Value = [double]"4.98285632E-25"

# Scientific-notation number format for the p-value row
$ws.Range("B8").NumberFormat = "0.00E+00"
$ws.Range("C8:D8").NumberFormat = "0.00E+00"
$ws.Range("F8").NumberFormat = "0.00E+00"
$ws.Range("I8:K8").NumberFormat = "0.00E+00"
$ws.Range("E8").NumberFormat = "0.00E+00"
$ws.Range("G8:H8").NumberFormat = "0.00E+00"
$ws.Range("L8").NumberFormat = "0.00E+00"

# Highlight the statistically-significant (near-zero) p-values
$ws.Range("E8").Interior.Color = $greenFill
$ws.Range("G8:H8").Interior.Color = $greenFill
$ws.Range("L8").Interior.Color = $greenFill

# --- Explanatory notes ---
$ws.Range("B11").Value = "Valores muy cercanos a cero indican alta relacion entre X e Y"
$ws.Range("B12").Value = "Tambien indican que es muy poco probable encontrar relacion entre X e Y por azar"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 24.17
$ws.Columns.Item(3).ColumnWidth = 20.42
$ws.Columns.Item(4).ColumnWidth = 23.48
$ws.Columns.Item(9).ColumnWidth = 18.77

# --- Final selection ---
$ws.Range("B13").Select()
